$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values per repulled data
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 1
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 1
$ws.Range("F14").Value = 3
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F29").Value = -5
$ws.Range("F36").Value = -2
$ws.Range("F37").Value = -5
$ws.Range("F38").Value = 3
$ws.Range("F40").Value = -2
$ws.Range("F43").Value = -1
$ws.Range("F44").Value = 1
